# Added code for advance TestNG
# The test data sheet (Sheet1) previously marked the vtiger login scenario's
# Status cell (E2) as the literal "123" (a placeholder / stray value). The
# updated test data expects it to read "pass" instead, and the workbook was
# left with the cursor/selection sitting on E2 (instead of D3) at a larger
# 160% zoom level, as last edited in the spreadsheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the Status value used by the vtiger test row.
$ws.Range("E2").Value = "pass"

# Reflect the author's last on-screen selection/zoom state for the sheet.
$ws.Activate()
$ws.Range("E2").Select()
$excel.ActiveWindow.Zoom = 160
